$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.535651683807373
$ws.Range("E2").Value = 2071.077284009994
$ws.Range("F2").Value = 0.09619130775575722
$ws.Range("G2").Value = 0.07223841429772591
$ws.Range("H2").Value = 0.06676356350911204
$ws.Range("I2").Value = 0.06079769944338075
$ws.Range("J2").Value = 0.05688322710620693
$ws.Range("K2").Value = 0.05201052365066593
$ws.Range("L2").Value = 0.04790697293061796
$ws.Range("M2").Value = 0.04700289203824763
$ws.Range("N2").Value = 0.04589928630816623
$ws.Range("O2").Value = 0.04474933810467559
$ws.Range("P2").Value = 0.04385171672781851
$ws.Range("Q2").Value = 0.04368924682026759
$ws.Range("R2").Value = 0.04273167645536627
$ws.Range("S2").Value = 0.04226820158176906
$ws.Range("T2").Value = 0.04181561970001633
$ws.Range("U2").Value = 0.04121798640775895
$ws.Range("V2").Value = 0.04099929559833988
$ws.Range("W2").Value = 0.04089451099853177
$ws.Range("X2").Value = 0.04037187688128643
$ws.Range("Y2").Value = 0.04037187688128643

$ws.Range("C3").Value = 0.5624661445617676
$ws.Range("E3").Value = 2057.554492135037
$ws.Range("F3").Value = 0.1040496186592716
$ws.Range("G3").Value = 0.07852545682778474
$ws.Range("H3").Value = 0.06799769908357159
$ws.Range("I3").Value = 0.05995013616928857
$ws.Range("J3").Value = 0.05436748848585393
$ws.Range("K3").Value = 0.05262684830478694
$ws.Range("L3").Value = 0.05040055608563304
$ws.Range("M3").Value = 0.04897208881178181
$ws.Range("N3").Value = 0.04706496158420085
$ws.Range("O3").Value = 0.04566894768263621
$ws.Range("P3").Value = 0.04529196877940888
$ws.Range("Q3").Value = 0.0450072194052421
$ws.Range("R3").Value = 0.04345335319870472
$ws.Range("S3").Value = 0.04220894050307804
$ws.Range("T3").Value = 0.04195759414312396
$ws.Range("U3").Value = 0.04077691356683637
$ws.Range("V3").Value = 0.04068566410494948
$ws.Range("W3").Value = 0.04067129025268212
$ws.Range("X3").Value = 0.040388671062328
$ws.Range("Y3").Value = 0.04010827470048806

$ws.Range("C4").Value = 0.6526522636413574
$ws.Range("E4").Value = 2045.657615153921
$ws.Range("F4").Value = 0.1058966226493461
$ws.Range("G4").Value = 0.08168233285313317
$ws.Range("H4").Value = 0.06520092173693254
$ws.Range("I4").Value = 0.06190800821996038
$ws.Range("J4").Value = 0.05388679412764582
$ws.Range("K4").Value = 0.05158533012751397
$ws.Range("L4").Value = 0.05043725681206127
$ws.Range("M4").Value = 0.04780562854854443
$ws.Range("N4").Value = 0.04594422824248828
$ws.Range("O4").Value = 0.04457799036002094
$ws.Range("P4").Value = 0.04391200502607547
$ws.Range("Q4").Value = 0.0432908698470944
$ws.Range("R4").Value = 0.04189930247207337
$ws.Range("S4").Value = 0.04138855878638478
$ws.Range("T4").Value = 0.04102833526821031
$ws.Range("U4").Value = 0.04102499667620754
$ws.Range("V4").Value = 0.04078021788252928
$ws.Range("W4").Value = 0.04017908764446591
$ws.Range("X4").Value = 0.04017908764446591
$ws.Range("Y4").Value = 0.03987636676713296

$ws.Range("C5").Value = 0.5312309265136719
$ws.Range("E5").Value = 2128.316071840802
$ws.Range("F5").Value = 0.1014106158508676
$ws.Range("G5").Value = 0.07734609197650608
$ws.Range("H5").Value = 0.0672092307138519
$ws.Range("I5").Value = 0.06034149114522559
$ws.Range("J5").Value = 0.05295198064463549
$ws.Range("K5").Value = 0.05164385609715772
$ws.Range("L5").Value = 0.05120934465906372
$ws.Range("M5").Value = 0.04749055353248344
$ws.Range("N5").Value = 0.04714585341636379
$ws.Range("O5").Value = 0.04589705072690506
$ws.Range("P5").Value = 0.04471134070215605
$ws.Range("Q5").Value = 0.04350719542978648
$ws.Range("R5").Value = 0.04350719542978648
$ws.Range("S5").Value = 0.04334102314932626
$ws.Range("T5").Value = 0.04334102314932626
$ws.Range("U5").Value = 0.04284520127094011
$ws.Range("V5").Value = 0.04238194029115867
$ws.Range("W5").Value = 0.04214751964972674
$ws.Range("X5").Value = 0.0417914160602707
$ws.Range("Y5").Value = 0.04148764272594155

$ws.Range("C6").Value = 0.5312235355377197
$ws.Range("E6").Value = 2089.482074480426
$ws.Range("F6").Value = 0.1070702801801101
$ws.Range("G6").Value = 0.08224452964030782
$ws.Range("H6").Value = 0.06704998608794363
$ws.Range("I6").Value = 0.06237069232982922
$ws.Range("J6").Value = 0.05663261033094853
$ws.Range("K6").Value = 0.05035931689133639
$ws.Range("L6").Value = 0.04974081658381365
$ws.Range("M6").Value = 0.04794766821577437
$ws.Range("N6").Value = 0.04726705088474455
$ws.Range("O6").Value = 0.04596727177957116
$ws.Range("P6").Value = 0.04465097466534084
$ws.Range("Q6").Value = 0.04327076948844024
$ws.Range("R6").Value = 0.04310315915729114
$ws.Range("S6").Value = 0.04238366874698696
$ws.Range("T6").Value = 0.04189124902669924
$ws.Range("U6").Value = 0.04156098884270207
$ws.Range("V6").Value = 0.0412698971866995
$ws.Range("W6").Value = 0.0412698971866995
$ws.Range("X6").Value = 0.0409294147286717
$ws.Range("Y6").Value = 0.04073064472671394

$ws.Range("C7").Value = 0.5468769073486328
$ws.Range("E7").Value = 2131.481860065461
$ws.Range("F7").Value = 0.101986988194821
$ws.Range("G7").Value = 0.07880670922953129
$ws.Range("H7").Value = 0.06808897355705008
$ws.Range("I7").Value = 0.05985535578245233
$ws.Range("J7").Value = 0.05771031829357398
$ws.Range("K7").Value = 0.05487149031878988
$ws.Range("L7").Value = 0.05088061676398413
$ws.Range("M7").Value = 0.04884595960345284
$ws.Range("N7").Value = 0.04873356011969549
$ws.Range("O7").Value = 0.04722432012245849
$ws.Range("P7").Value = 0.04597338384321248
$ws.Range("Q7").Value = 0.04452845603562638
$ws.Range("R7").Value = 0.04403729197979303
$ws.Range("S7").Value = 0.04365048283670805
$ws.Range("T7").Value = 0.04256797770551562
$ws.Range("U7").Value = 0.04256797770551562
$ws.Range("V7").Value = 0.04235309423509115
$ws.Range("W7").Value = 0.04191376068224167
$ws.Range("X7").Value = 0.04174204147940629
$ws.Range("Y7").Value = 0.04154935399737739

$ws.Range("C8").Value = 0.5469000339508057
$ws.Range("E8").Value = 2023.597818987548
$ws.Range("F8").Value = 0.10635682504802
$ws.Range("G8").Value = 0.0787425836271604
$ws.Range("H8").Value = 0.06317020909966578
$ws.Range("I8").Value = 0.06083868724956316
$ws.Range("J8").Value = 0.05224694444057304
$ws.Range("K8").Value = 0.05079948933075137
$ws.Range("L8").Value = 0.04679759354251271
$ws.Range("M8").Value = 0.04451570961916811
$ws.Range("N8").Value = 0.04260177248118992
$ws.Range("O8").Value = 0.04260177248118992
$ws.Range("P8").Value = 0.0424792610823708
$ws.Range("Q8").Value = 0.04185941795933858
$ws.Range("R8").Value = 0.04115732576758711
$ws.Range("S8").Value = 0.04058292982752929
$ws.Range("T8").Value = 0.04058277460324134
$ws.Range("U8").Value = 0.04012328508221332
$ws.Range("V8").Value = 0.03999567083629275
$ws.Range("W8").Value = 0.0397885404476403
$ws.Range("X8").Value = 0.03968658582795639
$ws.Range("Y8").Value = 0.03944635124732061

$ws.Range("C9").Value = 0.5625021457672119
$ws.Range("E9").Value = 2188.927713588515
$ws.Range("F9").Value = 0.09981227416678248
$ws.Range("G9").Value = 0.08105305313694013
$ws.Range("H9").Value = 0.06504013132723083
$ws.Range("I9").Value = 0.0642658463312468
$ws.Range("J9").Value = 0.0566330759569494
$ws.Range("K9").Value = 0.05484037832781292
$ws.Range("L9").Value = 0.05245023030293218
$ws.Range("M9").Value = 0.0506326366221254
$ws.Range("N9").Value = 0.05040430636392628
$ws.Range("O9").Value = 0.04912444697181451
$ws.Range("P9").Value = 0.04776848210919447
$ws.Range("Q9").Value = 0.04611261178391163
$ws.Range("R9").Value = 0.04512115879110108
$ws.Range("S9").Value = 0.04489355553977513
$ws.Range("T9").Value = 0.0436685956492556
$ws.Range("U9").Value = 0.0436685956492556
$ws.Range("V9").Value = 0.04349667922416686
$ws.Range("W9").Value = 0.04325422737153308
$ws.Range("X9").Value = 0.04304032999293108
$ws.Range("Y9").Value = 0.04266915621030242

$ws.Range("C10").Value = 0.5312457084655762
$ws.Range("E10").Value = 2165.246174192705
$ws.Range("F10").Value = 0.1015685571709321
$ws.Range("G10").Value = 0.07964939477719087
$ws.Range("H10").Value = 0.07349770844569661
$ws.Range("I10").Value = 0.0614359846432076
$ws.Range("J10").Value = 0.05682152990597198
$ws.Range("K10").Value = 0.0526554293413596
$ws.Range("L10").Value = 0.05081796689348621
$ws.Range("M10").Value = 0.04975104484224832
$ws.Range("N10").Value = 0.04811227984098332
$ws.Range("O10").Value = 0.04688742127765749
$ws.Range("P10").Value = 0.0464303792717419
$ws.Range("Q10").Value = 0.04548788139789175
$ws.Range("R10").Value = 0.04422384413117018
$ws.Range("S10").Value = 0.04399170218543218
$ws.Range("T10").Value = 0.04380502205492436
$ws.Range("U10").Value = 0.04333447037814719
$ws.Range("V10").Value = 0.04296702883607389
$ws.Range("W10").Value = 0.04268705812098521
$ws.Range("X10").Value = 0.04232799530331249
$ws.Range("Y10").Value = 0.04220752776204102

$ws.Range("C11").Value = 0.5625011920928955
$ws.Range("E11").Value = 2015.003410304233
$ws.Range("F11").Value = 0.1007093365244257
$ws.Range("G11").Value = 0.08042950830259647
$ws.Range("H11").Value = 0.07078516099194082
$ws.Range("I11").Value = 0.06317145840179626
$ws.Range("J11").Value = 0.05702202550025649
$ws.Range("K11").Value = 0.05361338590460223
$ws.Range("L11").Value = 0.05052797549419286
$ws.Range("M11").Value = 0.04772597507931552
$ws.Range("N11").Value = 0.04538390676528678
$ws.Range("O11").Value = 0.04400997132533806
$ws.Range("P11").Value = 0.04326183534608958
$ws.Range("Q11").Value = 0.04233703784903481
$ws.Range("R11").Value = 0.0417885173722594
$ws.Range("S11").Value = 0.04170332060492521
$ws.Range("T11").Value = 0.0411024583268873
$ws.Range("U11").Value = 0.04079794178216076
$ws.Range("V11").Value = 0.04008594606816723
$ws.Range("W11").Value = 0.03956030520877849
$ws.Range("X11").Value = 0.03956030520877849
$ws.Range("Y11").Value = 0.03927881891431251
